$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.370.50"
$ws.Range("E2").Value = "  +3.02%  "
$ws.Range("D3").Value = "2.096.56"
$ws.Range("E3").Value = "  +4.68%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.666"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +24.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("E12").Value = "  +8.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.30%  "
$ws.Range("D14").Value = "2.403.19"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.841"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.00%  "
$ws.Range("D16").Value = "2.108.96"
$ws.Range("E16").Value = "  +5.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.85%  "
$ws.Range("D18").Value = "37.298.35"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.97%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  +5.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.58%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.22%  "
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.54%  "
$ws.Range("E32").Value = "  +27.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0617"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.44%  "
$ws.Range("E35").Value = "  +10.77%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.13%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.90%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0227"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.09%  "
$ws.Range("E43").Value = "  +5.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.16%  "
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +101.68%  "
$ws.Range("D48").Value = "1.322.37"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  +6.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.05%  "
$ws.Range("E51").Value = "  +7.28%  "
